$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 176
$ws1.Range("F4").Value = 133
$ws1.Range("F5").Value = 1288
$ws1.Range("F6").Value = 17962
$ws1.Range("F7").Value = 354
$ws1.Range("F10").Value = 6776
$ws1.Range("F13").Value = 12
$ws1.Range("F14").Value = 107
$ws1.Range("F19").Value = 206
$ws1.Range("F22").Value = 35
$ws1.Range("F26").Value = 975
$ws1.Range("F30").Value = 15
$ws1.Range("F31").Value = 9
$ws1.Range("F32").Value = 68
$ws1.Range("F33").Value = 11994
$ws1.Range("F34").Value = 1273
$ws1.Range("F37").Value = 267

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 176
$ws4.Range("F4").Value = 133
$ws4.Range("F5").Value = 1288
$ws4.Range("F6").Value = 17962
$ws4.Range("F7").Value = 354
$ws4.Range("F10").Value = 6776
$ws4.Range("F13").Value = 12
$ws4.Range("F14").Value = 107
$ws4.Range("F19").Value = 206
$ws4.Range("F22").Value = 35
$ws4.Range("F26").Value = 975
$ws4.Range("F32").Value = 15
$ws4.Range("F33").Value = 9
$ws4.Range("F34").Value = 68
$ws4.Range("F35").Value = 11994
$ws4.Range("F36").Value = 1273
$ws4.Range("F39").Value = 267
